$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 2-40: price (D) and volume (E) columns ---
$ws.Range("D2").Value = '25.556.04'
$ws.Range("E2").Value = '  +1.79%  '
$ws.Range("D3").Value = '1.664.36'
$ws.Range("E3").Value = '  +0.83%  '
$ws.Range("D4").Value = '''0.9991'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''235.87'
$ws.Range("E5").Value = '  -0.64%  '
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").Value = '''0.4795'
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '''0.2613'
$ws.Range("E8").Value = '  -0.53%  '
$ws.Range("D9").Value = '''0.06148'
$ws.Range("E9").Value = '  +1.87%  '
$ws.Range("D10").Value = '''0.07089'
$ws.Range("E10").Value = '  -0.10%  '
$ws.Range("D11").Value = '1.673.09'
$ws.Range("E11").Value = '  +1.14%  '
$ws.Range("D12").Value = '''14.74'
$ws.Range("E12").Value = '  +1.76%  '
$ws.Range("D13").Value = '''0.5905'
$ws.Range("E13").Value = '  -4.64%  '
$ws.Range("D14").Value = '''4.367'
$ws.Range("E14").Value = '  -4.37%  '
$ws.Range("D15").Value = '''74.40'
$ws.Range("E15").Value = '  +1.69%  '
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").Value = '''0.9998'
$ws.Range("E17").Value = '  +0.04%  '
$ws.Range("D18").Value = '25.538.66'
$ws.Range("E18").Value = '  +1.82%  '
$ws.Range("D19").Value = '''0.000006748'
$ws.Range("E19").Value = '  +2.96%  '
$ws.Range("E20").Value = '  +0.27%  '
$ws.Range("D21").Value = '1.878.51'
$ws.Range("E21").Value = '  +1.04%  '
$ws.Range("D22").Value = '''4.422'
$ws.Range("E22").Value = '  +0.36%  '
$ws.Range("D23").Value = '''8.652'
$ws.Range("E23").Value = '  +2.16%  '
$ws.Range("D24").Value = '''5.299'
$ws.Range("E24").Value = '  +1.20%  '
$ws.Range("D25").Value = '''134.38'
$ws.Range("E25").Value = '  +0.40%  '
$ws.Range("E26").Value = '  +1.94%  '
$ws.Range("D27").Value = '''1.397'
$ws.Range("E27").Value = '  +0.17%  '
$ws.Range("D28").Value = '''104.62'
$ws.Range("E28").Value = '  +2.72%  '
$ws.Range("D29").Value = '''1.684'
$ws.Range("E29").Value = '  -0.38%  '
$ws.Range("D30").Value = '''3.956'
$ws.Range("E30").Value = '  +4.46%  '
$ws.Range("E31").Value = '  +2.79%  '
$ws.Range("D32").Value = '''0.07622'
$ws.Range("E32").Value = '  -3.76%  '
$ws.Range("D33").Value = '''0.9994'
$ws.Range("E33").Value = '  +0.03%  '
$ws.Range("D34").Value = '''0.04309'
$ws.Range("E34").Value = '  -5.06%  '
$ws.Range("D35").Value = '''2.617'
$ws.Range("E35").Value = '  +0.46%  '
$ws.Range("D36").Value = '''0.6122'
$ws.Range("E36").Value = '  +5.75%  '
$ws.Range("D37").Value = '''0.9480'
$ws.Range("E37").Value = '  +0.55%  '
$ws.Range("D38").Value = '''2.608'
$ws.Range("E38").Value = '  -0.43%  '
$ws.Range("D39").Value = '''0.8472'
$ws.Range("E39").Value = '  +1.13%  '
$ws.Range("E40").Value = '  +0.02%  '

# --- Insert new row 41 (PaxosStandard) and shift rows 41-51 down to 42-52,
#     then drop the row that falls off the bottom (formerly USDD) to keep 51 rows ---
$ws.Rows.Item(41).Insert()
$ws.Rows.Item(52).Delete()

# Copy formatting for the new index cell A41 from a neighboring data row, then set its value
$ws.Range("A40").Copy()
$ws.Range("A41").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Column A is a plain sequential row index (0..49) and is untouched by the
# coin reshuffle, so restore it to a clean 0-based sequence for every data row.
for ($i = 0; $i -le 49; $i++) {
    $ws.Range("A" + ($i + 2)).Value = $i
}

# --- Set final values for rows 41-51 (Coin, Link, Price, Volume) ---
# Row 41
$ws.Range("B41").Value = 'PaxosStandard'
$ws.Range("C41").Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
$ws.Range("D41").Value = '''0.9999'
$ws.Range("E41").Value = '  -0.05%  '
# Row 42
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '''0.01494'
$ws.Range("E42").Value = '  -2.90%  '
# Row 43
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").Value = '''1.872'
$ws.Range("E43").Value = '  +2.96%  '
# Row 44
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '''97.78'
$ws.Range("E44").Value = '  -0.82%  '
# Row 45
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '''0.3756'
$ws.Range("E45").Value = '  +1.41%  '
# Row 46
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").Value = '''4.691'
$ws.Range("E46").Value = '  -2.08%  '
# Row 47
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").Value = '''0.1119'
$ws.Range("E47").Value = '  -0.99%  '
# Row 48
$ws.Range("B48").Value = 'Aptos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D48").Value = '''6.209'
$ws.Range("E48").Value = '  +2.93%  '
# Row 49
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '''0.05259'
$ws.Range("E49").Value = '  +1.93%  '
# Row 50
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").Value = '''29.35'
$ws.Range("E50").Value = '  -0.80%  '
# Row 51
$ws.Range("B51").Value = 'TrueUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range("D51").Value = '''1.001'
$ws.Range("E51").Value = '  +0.16%  '
